$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "ops" header relabeling: the assignment columns move from the old
# Emp-ID / QA wording to the new Typist / Typist QC wording.
$ws.Range("C1").Value = "Typist"
$ws.Range("D1").Value = "Typist QC"

# Sample data row: updated employee codes and corrected status text
# (WIP -> Typing) for the order-form/order-page example row.
$ws.Range("C2").Value = "SIPL0102"
$ws.Range("D2").Value = "SIPL5317"
$ws.Range("L2").Value = "Typing"

# Leave the cursor where it was left in the saved workbook.
$null = $ws.Range("E5").Select()
